$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PY")

# Insert 2 new rows before row 27, shifting the rest of the table down.
$ws.Range("A27:A28").EntireRow.Insert()

# Fill the three new rows (127/128/129 order in the shared-string table
# comes from the order the strings were first typed).
$ws.Range("B28").Value = "cmd - pip freeze (xem những thư viện đã cài)"
$ws.Range("B26").Value = "cmd - pip install … (… tên thư viện muốn cài)"
$ws.Range("B27").Value = "cmd - pip install -r … (… file chứa danh sách thư viện cần cài - vd: requirements.txt)"

# Reflect the author's on-screen viewport/selection when the edit was saved.
$wsGit = $wb.Worksheets.Item("GIT")
$wsGit.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 8
$wsGit.Range("C27").Select() | Out-Null

$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 15
$ws.Range("C24").Select() | Out-Null
$wb.Windows.Item(1).ActiveSheet = $ws
